$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Gimmick sheet: add a new "damage" (int) column in F, and a new gimmick
# row (gim_damage_wood_step) so gimmicks can deal damage to characters.
# ---------------------------------------------------------------------------
$gim = $wb.Worksheets.Item("Gimmick")

# Column F formatting mirrors column D (int-typed stat column).
$gim.Range("D1:D4").Copy($gim.Range("F1:F4"))
$gim.Range("F1").Value = "damage"

$gim.Range("D5").Copy($gim.Range("F5"))
$gim.Range("F5").Value = 0

$gim.Range("A5:E5").Copy($gim.Range("A6:E6"))
$gim.Range("F5").Copy($gim.Range("F6"))
$gim.Range("A6").Value = 2
$gim.Range("B6").Value = "gim_damage_wood_step"
$gim.Range("C6").Value = 1
$gim.Range("D6").Value = 5
$gim.Range("E6").Value = 5
$gim.Range("F6").Value = 1

# Column widths: B widens to fit the longer asset name, C:F share the
# standard stat-column width (stored width unit is ColumnWidth + 5/6,
# quantized to the nearest 1/6 of a character by the engine).
$gim.Columns.Item(2).ColumnWidth = 24 - (5/6)
$gim.Columns.Item(6).ColumnWidth = 13

# ---------------------------------------------------------------------------
# Character sheet: add a new "hp" (int) column in H, matching the existing
# fieldName/int header pattern used by the other stat columns.
# ---------------------------------------------------------------------------
$chr = $wb.Worksheets.Item("Character")

# Copy the formatting of the neighbouring column G (header/blank/footer/data
# styles) into column H, then overwrite the values that differ.
$chr.Range("G1:G5").Copy($chr.Range("H1:H5"))

$chr.Range("H1").Value = "hp"
$chr.Range("H5").Value = 3

# Column H should use the same 14.5 width as columns C:G (stored width unit
# is ColumnWidth + 5/6).
$chr.Columns.Item(8).ColumnWidth = 14.5 - (5/6)

$chr.Activate()
$chr.Range("E30:E31").Select()

# Leave the Gimmick sheet as the active tab (it was already the active tab
# before this edit) and restore its own selection.
$gim.Activate()
$gim.Range("D15").Select()
